# Selenium second commit and push
# Rework the "data" sheet from a product/price list into a
# username/password style credential list with hyperlinked emails,
# keeping the iphone/redmi/nokia price rows intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Insert a new column A. This pushes the old Product Name / Product
# Price columns (and their auto-fit "bestFit" widths) one slot to the
# right, into B and C.
$ws.Columns.Item(1).Insert()

# The iphone / redmi / nokia rows are not edited - just slid back from
# B/C into A/B so the price table keeps working.
foreach ($r in 4, 6, 7) {
    $name = $ws.Range("B$r").Value2
    $price = $ws.Range("C$r").Value2
    $ws.Range("A$r").Value = $name
    $ws.Range("B$r").Value = $price
    $ws.Range("C$r").ClearContents()
}

# --- Row 1: header ---
$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "passward"
$ws.Range("C1").ClearContents()

# --- Row 2: samsung login ---
$ws.Range("B2").Value = "biswa234"
$ws.Range("C2").ClearContents()
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:samsung234@mail.com", "", "", "samsung234@mail.com") | Out-Null

# --- Row 3: vivo login ---
$ws.Range("B3").Value = "viv023444"
$ws.Range("C3").ClearContents()
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:vivo2343@gmail.com", "", "", "vivo2343@gmail.com") | Out-Null

# --- Row 5: Nothing login ---
$ws.Range("B5").Value = 40000
$ws.Range("C5").ClearContents()
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:Nothing2323@gmail.com", "", "", "Nothing2323@gmail.com") | Out-Null

# --- Row 8: extra password only ---
$ws.Range("B8").Value = "bhohf345"

# --- Row 9: balia login ---
$ws.Range("B9").Value = "bhibu567"
$ws.Hyperlinks.Add($ws.Range("A9"), "mailto:balia456@gmail.com", "", "", "balia456@gmail.com") | Out-Null

# --- Row 10: phone numbers ---
$ws.Range("A10").Value = 8260539737
$ws.Range("B10").Value = 9861000762

# Widen the new first column for the long e-mail addresses it now holds
# (best attainable approximation of the 20.5546875 bestFit width Excel
# itself would compute for "Nothing2323@gmail.com" / "samsung234@mail.com").
$ws.Columns.Item(1).ColumnWidth = 19.6

# Leave the selection where the author's session ended up.
$ws.Range("B10").Select() | Out-Null

Write-Output "done"
